# Updates Leve-profit market-data columns (H:N) on several worksheets
# to refresh currentAveragePrice* / LevePrice* / LeveProfit* figures
# pulled in by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3576498.8
$ws.Range("J17").Value = 4005548.2
$ws.Range("L17").Value = 12016644.6
$ws.Range("N17").Value = -12016980.6
$ws.Range("H116").Value = 4881.1333
$ws.Range("I116").Value = 2435
$ws.Range("J116").Value = 5492.6665
$ws.Range("K116").Value = 2435
$ws.Range("L116").Value = 5492.6665
$ws.Range("M116").Value = 1007
$ws.Range("N116").Value = -12376.6665
$ws.Range("H132").Value = 3488.7778
$ws.Range("I132").Value = 3686.5652
$ws.Range("J132").Value = 2351.5
$ws.Range("K132").Value = 11059.6956
$ws.Range("L132").Value = 7054.5
$ws.Range("M132").Value = -8529.695599999999
$ws.Range("N132").Value = -12114.5
$ws.Range("H137").Value = 1796.1666
$ws.Range("I137").Value = 1657.2858
$ws.Range("J137").Value = 1884.5454
$ws.Range("K137").Value = 4971.857400000001
$ws.Range("L137").Value = 5653.6362
$ws.Range("M137").Value = -2421.857400000001
$ws.Range("N137").Value = -10753.6362
$ws.Range("H138").Value = 10754940
$ws.Range("I138").Value = 22728102
$ws.Range("J138").Value = 3528.2856
$ws.Range("K138").Value = 68184306
$ws.Range("L138").Value = 10584.8568
$ws.Range("M138").Value = -68179166
$ws.Range("N138").Value = -20864.8568

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3039.25
$ws.Range("I45").Value = 2906.3809
$ws.Range("J45").Value = 3292.9092
$ws.Range("K45").Value = 2906.3809
$ws.Range("L45").Value = 3292.9092
$ws.Range("M45").Value = -2529.3809
$ws.Range("N45").Value = -4046.9092
$ws.Range("H61").Value = 621995.9
$ws.Range("I61").Value = 751383.5600000001
$ws.Range("J61").Value = 934.8
$ws.Range("K61").Value = 751383.5600000001
$ws.Range("L61").Value = 934.8
$ws.Range("M61").Value = -751171.5600000001
$ws.Range("N61").Value = -1358.8
$ws.Range("H74").Value = 38463828
$ws.Range("I74").Value = 38463828
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 38463828
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -38462954
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 38463828
$ws.Range("I77").Value = 38463828
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 192319140
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -192314772
$ws.Range("N77").ClearContents()
$ws.Range("H136").Value = 621995.9
$ws.Range("I136").Value = 751383.5600000001
$ws.Range("J136").Value = 934.8
$ws.Range("K136").Value = 2254150.68
$ws.Range("L136").Value = 2804.4
$ws.Range("M136").Value = -2251600.68
$ws.Range("N136").Value = -7904.4

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3752.8108
$ws.Range("I134").Value = 4081.1724
$ws.Range("K134").Value = 12243.5172
$ws.Range("M134").Value = -9708.5172

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3123.5952
$ws.Range("I31").Value = 1792.3448
$ws.Range("J31").Value = 6093.3076
$ws.Range("K31").Value = 1792.3448
$ws.Range("L31").Value = 6093.3076
$ws.Range("M31").Value = -1497.3448
$ws.Range("N31").Value = -6683.3076
$ws.Range("H34").Value = 3123.5952
$ws.Range("I34").Value = 1792.3448
$ws.Range("J34").Value = 6093.3076
$ws.Range("K34").Value = 1792.3448
$ws.Range("L34").Value = 6093.3076
$ws.Range("M34").Value = -1590.3448
$ws.Range("N34").Value = -6497.3076
$ws.Range("H132").Value = 2339.0303
$ws.Range("I132").Value = 1755.4375
$ws.Range("K132").Value = 5266.3125
$ws.Range("M132").Value = -2736.3125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8245.571
$ws.Range("I80").Value = 1107.6
$ws.Range("J80").Value = 12211.111
$ws.Range("K80").Value = 3322.8
$ws.Range("L80").Value = 36633.333
$ws.Range("M80").Value = -2386.8
$ws.Range("N80").Value = -38505.333
$ws.Range("H83").Value = 8245.571
$ws.Range("I83").Value = 1107.6
$ws.Range("J83").Value = 12211.111
$ws.Range("K83").Value = 9968.4
$ws.Range("L83").Value = 109899.999
$ws.Range("M83").Value = -5288.4
$ws.Range("N83").Value = -119259.999
$ws.Range("H118").Value = 45459956
$ws.Range("J118").Value = 9732.833000000001
$ws.Range("L118").Value = 29198.499
$ws.Range("N118").Value = -31684.499
$ws.Range("H119").Value = 4810
$ws.Range("I119").Value = 3400
$ws.Range("K119").Value = 10200
$ws.Range("M119").Value = -5362
$ws.Range("H131").Value = 695.6598
$ws.Range("J131").Value = 716.8068
$ws.Range("L131").Value = 2150.4204
$ws.Range("N131").Value = -12230.4204
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 43366.31
$ws.Range("I132").Value = 5313.5
$ws.Range("K132").Value = 15940.5
$ws.Range("M132").Value = -13410.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47624372
$ws.Range("I7").Value = 83336510
$ws.Range("K7").Value = 83336510
$ws.Range("M7").Value = -83336398
$ws.Range("H68").Value = 2999.25
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 2999.25
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -22483
$ws.Range("H126").Value = 47624372
$ws.Range("I126").Value = 83336510
$ws.Range("K126").Value = 250009530
$ws.Range("M126").Value = -250007060
$ws.Range("H132").Value = 287713.84
$ws.Range("I132").Value = 294711.72
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 884135.1599999999
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -881605.1599999999
$ws.Range("N132").Value = -7460
$ws.Range("H136").Value = 931.05
$ws.Range("I136").Value = 825.75
$ws.Range("K136").Value = 2477.25
$ws.Range("M136").Value = 72.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1423078.6
$ws.Range("I113").Value = 782.7857
$ws.Range("K113").Value = 2348.3571
$ws.Range("M113").Value = -178.3571000000002
$ws.Range("H122").Value = 1300.1666
$ws.Range("I122").Value = 1313
$ws.Range("K122").Value = 3939
$ws.Range("M122").Value = -1489
$ws.Range("H136").Value = 15882756
$ws.Range("I136").Value = 22441262
$ws.Range("J136").Value = 4267.5264
$ws.Range("K136").Value = 67323786
$ws.Range("L136").Value = 12802.5792
$ws.Range("M136").Value = -67321236
$ws.Range("N136").Value = -17902.5792
